$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 4387.2, 4741, 3890, 0.05120588938395183),
    @(1, 4189.9, 4512, 3577, 0.05578962961832682),
    @(2, 4850.6, 5145, 4607, 0.05418267250061035),
    @(3, 5112.033333333334, 5522, 4570, 0.054460875193278),
    @(4, 4328.2, 4562, 3976, 0.05216947396596273),
    @(5, 4986, 5314, 4421, 0.05562014579772949),
    @(6, 4425.866666666667, 4815, 3915, 0.0545870304107666),
    @(7, 4518.7, 4878, 4143, 0.0554798444112142),
    @(8, 4680.5, 5155, 4096, 0.05877664089202881),
    @(9, 4780.966666666666, 5243, 4389, 0.05096296469370524)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
    $ws.Cells.Item($row, 5).Value = $data[$i][4]
}
